# Edit script for "12 Step 04 - Understanding Camel Terminology and Architecture.docx"
#
# Summary of the changes applied (see commit diff):
#   1. The paragraph "Camel Context:" (the bold heading that introduces the
#      "Camel Context = Routes + Components ..." bullet) is turned into a new
#      "Jatin: " heading, followed by a brand-new bullet "Route: A sequence of
#      processing steps being performed on msg when msg travels from source to
#      destination." and then the original "Camel Context:" heading is
#      re-created right after it (so the net effect is: two new paragraphs are
#      inserted in front of the untouched "Camel Context:" heading).
#   2. The <w:lastRenderedPageBreak/> marker moves from the "Route Processor:"
#      run to the "Filter Processor:" run (pure layout bookkeeping, caused by
#      the extra content added earlier in the document).
#
# Each block below is applied by replacing the OOXML of a precisely-targeted
# Range via Range.InsertXML(), which guarantees exact run/paragraph structure
# (no stray inherited rPr/bCs artifacts that plain Range.Text edits would
# leave behind).

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Camel Context:" -> "Jatin:" + new "Route: ..." bullet + restored
#    "Camel Context:" heading.
# ---------------------------------------------------------------------------
$camelContextHeading = Get-ParagraphByText $d "Camel Context:"

$jatinRouteCamelXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Jatin: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Route: </w:t></w:r><w:r><w:t xml:space="preserve">A sequence of </w:t></w:r><w:r><w:t xml:space="preserve">processing </w:t></w:r><w:r><w:t>steps being performed on msg when msg travels from source to destination.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Camel Context</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$camelContextHeading.Range.InsertXML($jatinRouteCamelXml)

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from "Route Processor:" to
#    "Filter Processor:".
# ---------------------------------------------------------------------------
$filterProcessorHeading = Get-ParagraphByText $d "Filter Processor:"
$routeProcessorHeading  = Get-ParagraphByText $d "Route Processor:"

$filterProcessorXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Filter Processor</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$routeProcessorXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Route Processor</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$filterProcessorHeading.Range.InsertXML($filterProcessorXml)
$routeProcessorHeading.Range.InsertXML($routeProcessorXml)
